$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'330.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.86%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'41.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'2.64%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.715"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.45%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08093"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.90%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.037"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'5.15%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'8.725"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.15%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'4.512"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-1.33%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'-0.68%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9217"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-2.33%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1251"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-0.03%"
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'-0.84%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'8.329"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-5.90%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09313"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'1.33%"
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'2.12%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.1055"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'9.61%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.001305"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.01%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.006138"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-6.44%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.381"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.30%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'-1.24%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1417"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-1.21%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2652"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'9.77%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04430"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.70%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001261"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.07%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004348"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.66%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'8.54%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D39").Value = "'0.02807"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'16.15%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05476"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'3.74%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007587"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.96%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.009968"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'14.19%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1422"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'0.28%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002119"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'0.64%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.01177"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'22.62%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006765"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-1.81%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.37%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.002979"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-5.55%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.002281"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'60.09%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'-0.37%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'-0.37%"
$ws.Range("E51").Style = "Normal"
